# Generate Report for Handoff
# Renames the two tracked files, refreshes their status/handoff metadata,
# and clears out the "Latest Target File" / "Latest Handback File" /
# hyperlink info that no longer applies.

$wb = $excel.ActiveWorkbook

$oldFile1 = "2c5e234e-fb30-4425-885b-b108c390ff7c"
$newFile1 = "73a0d2a1-89a1-481a-ab4c-4ce0e12a45c0"
$oldFile2 = "abf2a8c3-5959-415f-906b-fab2aa37853a"
$newFile2 = "ffff341494fd-1949-4346-98a4-80679563b85d"

$newStatus = "Ready for handoff"
$newHoDate = "2016-09-05 23:16:38"
$newHoDateZhCn = "2016-09-05 23:16:32"
$newHoDateDeDe = "2016-09-05 23:16:38"
$zeroDate = "0001-01-01 00:00:00"

$newZhCnHandoff1 = $newFile1 + ".5699311b2f211a1627904973773dd11918028c8f.zh-cn.xlf"
$newDeDeHandoff1 = $newFile1 + ".5699311b2f211a1627904973773dd11918028c8f.de-de.xlf"
$newZhCnHandoff2 = $newFile1 + ".5699311b2f211a1627904973773dd11918028c8f.zh-cn.xlf"
$newDeDeHandoff2 = $newFile1 + ".5699311b2f211a1627904973773dd11918028c8f.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(2,1).Value = "'" + $newFile1 + ".md"
$wsOverview.Cells.Item(2,2).Value = "'e2e\" + $newFile1 + ".md"
$wsOverview.Cells.Item(2,5).Value = $newStatus
$wsOverview.Cells.Item(2,6).Value = $newStatus
$wsOverview.Cells.Item(2,7).Value = $newHoDate

$wsOverview.Cells.Item(3,1).Value = "'" + $newFile2 + ".md"
$wsOverview.Cells.Item(3,2).Value = "'e2e\" + $newFile2 + ".md"
$wsOverview.Cells.Item(3,5).Value = $newStatus
$wsOverview.Cells.Item(3,6).Value = $newStatus
$wsOverview.Cells.Item(3,7).Value = $newHoDate

# Recreate the two hyperlinks (same targets as before, just refreshed display text)
$ovB2Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/" + $oldFile1 + ".md"
$ovB3Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/" + $oldFile2 + ".md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(2,2), $ovB2Addr, "", "", "e2e\" + $newFile1 + ".md")
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3,2), $ovB3Addr, "", "", "e2e\" + $newFile2 + ".md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(2,1).Value = "'" + $newFile1 + ".md"
$wsZhCn.Cells.Item(2,3).Value = $newStatus
$wsZhCn.Cells.Item(2,6).Value = "'False"
$wsZhCn.Cells.Item(2,6).Style = "Normal"
$wsZhCn.Cells.Item(2,7).Value = $newZhCnHandoff1
$wsZhCn.Cells.Item(2,8).Value = $newHoDateZhCn
$wsZhCn.Cells.Item(2,9).Value = ""
$wsZhCn.Cells.Item(2,9).Style = "Normal"
$wsZhCn.Cells.Item(2,10).Value = ""
$wsZhCn.Cells.Item(2,11).Value = $zeroDate

$wsZhCn.Cells.Item(3,1).Value = "'" + $newFile2 + ".md"
$wsZhCn.Cells.Item(3,3).Value = $newStatus
$wsZhCn.Cells.Item(3,6).Value = "'True"
$wsZhCn.Cells.Item(3,6).Style = "Normal"
$wsZhCn.Cells.Item(3,7).Value = $newZhCnHandoff2
$wsZhCn.Cells.Item(3,8).Value = $newHoDateZhCn
$wsZhCn.Cells.Item(3,9).Value = ""
$wsZhCn.Cells.Item(3,9).Style = "Normal"
$wsZhCn.Cells.Item(3,10).Value = ""
$wsZhCn.Cells.Item(3,11).Value = $zeroDate

$zhA2Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/" + $oldFile1 + ".md"
$zhA3Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/" + $oldFile2 + ".md"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(2,1), $zhA2Addr, "", "", $newFile1 + ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(3,1), $zhA3Addr, "", "", $newFile2 + ".md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(2,1).Value = "'" + $newFile1 + ".md"
$wsDeDe.Cells.Item(2,3).Value = $newStatus
$wsDeDe.Cells.Item(2,6).Value = "'False"
$wsDeDe.Cells.Item(2,6).Style = "Normal"
$wsDeDe.Cells.Item(2,7).Value = $newDeDeHandoff1
$wsDeDe.Cells.Item(2,8).Value = $newHoDateDeDe
$wsDeDe.Cells.Item(2,9).Value = ""
$wsDeDe.Cells.Item(2,9).Style = "Normal"
$wsDeDe.Cells.Item(2,10).Value = ""
$wsDeDe.Cells.Item(2,11).Value = $zeroDate

$wsDeDe.Cells.Item(3,1).Value = "'" + $newFile2 + ".md"
$wsDeDe.Cells.Item(3,3).Value = $newStatus
$wsDeDe.Cells.Item(3,6).Value = "'True"
$wsDeDe.Cells.Item(3,6).Style = "Normal"
$wsDeDe.Cells.Item(3,7).Value = $newDeDeHandoff2
$wsDeDe.Cells.Item(3,8).Value = $newHoDateDeDe
$wsDeDe.Cells.Item(3,9).Value = ""
$wsDeDe.Cells.Item(3,9).Style = "Normal"
$wsDeDe.Cells.Item(3,10).Value = ""
$wsDeDe.Cells.Item(3,11).Value = $zeroDate

$deA2Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/" + $oldFile1 + ".md"
$deA3Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02ddc6d51716c6ceb41d4a56fd86f4b2c84289d4/e2e/" + $oldFile2 + ".md"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(2,1), $deA2Addr, "", "", $newFile1 + ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(3,1), $deA3Addr, "", "", $newFile2 + ".md")

# ---------------------------------------------------------------------------
# Column widths: content got shorter, so shrink the affected columns to fit.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()

$wsZhCn.Columns.Item(3).AutoFit()
$wsZhCn.Columns.Item(9).AutoFit()
$wsZhCn.Columns.Item(10).AutoFit()

$wsDeDe.Columns.Item(3).AutoFit()
$wsDeDe.Columns.Item(9).AutoFit()
$wsDeDe.Columns.Item(10).AutoFit()
